$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (date / time number formats) of the last existing
# data row (28) down across the five new rows (29:33) so the new cells
# reuse the workbook's existing styles instead of minting new ones.
$ws.Range("A28:C28").Copy()
$ws.Range("A29:C33").PasteSpecial(-4122)

$data = @(
    @(43331, "Beastlord", 0.017905092592592594),
    @(43336, "PCAutoMechanic", 0.036331018518518519),
    @(43338, "PCAutoMechanic", 0.022349537037037032),
    @(43343, "ShadowOfWar", 0.036724537037037035),
    @(43346, "ShadowOfWar", 0.013530092592592594)
)

$startRow = 29
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $entry = $data[$i]

    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
}

$ws.Range("A34").Select()
